# "Additional test cases and test template"
#
# The credentials.xlsx fixture gets a couple of its sample rows reworked:
#   - row 6 (marcin.koziol@test.com / letmein123) is now flagged as a
#     valid set of credentials ("Yes")
#   - row 7 (bogdan.romanski@test.com / counttothree123) keeps its
#     existing "No" / invalid status
#   - row 8's test account is swapped for a new one
#     (marin.kozieradka@test.com / koziol333) and is also flagged valid
#   - the cursor is left parked on D15, as in the saved template

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 — rewrite the password (same text) and flip the verdict to Yes.
$ws.Range("B6").Value = "letmein123"
$ws.Range("C6").Value = "Yes"

# Row 7 — rewrite the password (same text); verdict stays No.
$ws.Range("B7").Value = "counttothree123"

# Row 8 — brand-new test account, flagged valid.
$ws.Range("A8").Value = "marin.kozieradka@test.com"
$ws.Range("B8").Value = "koziol333"
$ws.Range("C8").Value = "Yes"

# Leave the selection where the author left it when saving.
$ws.Range("D15").Select() | Out-Null
